$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 465.44446
$ws.Range("I4").Value = 423.625
$ws.Range("K4").Value = 423.625
$ws.Range("M4").Value = -309.625

$ws.Range("H15").Value = 258.25
$ws.Range("I15").Value = 258.25
$ws.Range("K15").Value = 774.75
$ws.Range("M15").Value = -605.75

$ws.Range("H17").Value = 5943301.5
$ws.Range("J17").Value = 5943301.5
$ws.Range("L17").Value = 17829904.5
$ws.Range("N17").Value = -17830240.5

$ws.Range("H28").Value = 199
$ws.Range("I28").Value = 107.5
$ws.Range("J28").Value = 382
$ws.Range("K28").Value = 107.5
$ws.Range("L28").Value = 382
$ws.Range("M28").Value = 377.5
$ws.Range("N28").Value = -1352

$ws.Range("H53").Value = 45796.816
$ws.Range("I53").Value = 132.44444
$ws.Range("J53").Value = 77410.62
$ws.Range("K53").Value = 132.44444
$ws.Range("L53").Value = 77410.62
$ws.Range("M53").Value = 504.55556
$ws.Range("N53").Value = -78684.62

$ws.Range("H86").Value = 3032.7778
$ws.Range("I86").Value = 3160.7856
$ws.Range("J86").Value = 2584.75
$ws.Range("K86").Value = 3160.7856
$ws.Range("L86").Value = 2584.75
$ws.Range("M86").Value = -2037.7856
$ws.Range("N86").Value = -4830.75

$ws.Range("H89").Value = 3032.7778
$ws.Range("I89").Value = 3160.7856
$ws.Range("J89").Value = 2584.75
$ws.Range("K89").Value = 15803.928
$ws.Range("L89").Value = 12923.75
$ws.Range("M89").Value = -10187.928
$ws.Range("N89").Value = -24155.75

$ws.Range("H112").Value = 2464056.5
$ws.Range("J112").Value = 2552044.2
$ws.Range("L112").Value = 7656132.600000001
$ws.Range("N112").Value = -7658348.600000001

$ws.Range("H121").Value = 1504.6
$ws.Range("I121").Value = 800
$ws.Range("J121").Value = 1541.6842
$ws.Range("K121").Value = 2400
$ws.Range("L121").Value = 4625.0526
$ws.Range("M121").Value = -653
$ws.Range("N121").Value = -8119.0526

$ws.Range("H129").Value = 1079.4756
$ws.Range("I129").Value = 837.5
$ws.Range("J129").Value = 1091.8846
$ws.Range("K129").Value = 2512.5
$ws.Range("L129").Value = 3275.6538
$ws.Range("M129").Value = 2487.5
$ws.Range("N129").Value = -13275.6538

$ws.Range("H132").Value = 1814.1515
$ws.Range("I132").Value = 1655.5938
$ws.Range("K132").Value = 4966.7814
$ws.Range("M132").Value = -2436.7814

$ws.Range("H137").Value = 1672.2128
$ws.Range("I137").Value = 1302.6285
$ws.Range("K137").Value = 3907.8855
$ws.Range("M137").Value = -1357.8855

$ws.Range("H141").Value = 3362.9048
$ws.Range("I141").Value = 2338.182
$ws.Range("J141").Value = 4490.1
$ws.Range("K141").Value = 7014.545999999999
$ws.Range("L141").Value = 13470.3
$ws.Range("M141").Value = -1834.545999999999
$ws.Range("N141").Value = -23830.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3141.22
$ws.Range("I32").Value = 3022.923
$ws.Range("J32").Value = 4337.3335
$ws.Range("K32").Value = 3022.923
$ws.Range("L32").Value = 4337.3335
$ws.Range("M32").Value = -2735.923
$ws.Range("N32").Value = -4911.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1265.8334
$ws.Range("J11").Value = 2847.5
$ws.Range("L11").Value = 2847.5
$ws.Range("N11").Value = -3127.5

$ws.Range("H20").Value = 26862
$ws.Range("I20").Value = 2160
$ws.Range("J20").Value = 84500
$ws.Range("K20").Value = 2160
$ws.Range("L20").Value = 84500
$ws.Range("M20").Value = -1913
$ws.Range("N20").Value = -84994

$ws.Range("H80").Value = 279.3889
$ws.Range("I80").Value = 90.5
$ws.Range("J80").Value = 333.35715
$ws.Range("K80").Value = 90.5
$ws.Range("L80").Value = 333.35715
$ws.Range("M80").Value = 907.5
$ws.Range("N80").Value = -2329.35715

$ws.Range("H83").Value = 279.3889
$ws.Range("I83").Value = 90.5
$ws.Range("J83").Value = 333.35715
$ws.Range("K83").Value = 452.5
$ws.Range("L83").Value = 1666.78575
$ws.Range("M83").Value = 4539.5
$ws.Range("N83").Value = -11650.78575

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2190.6775
$ws.Range("I132").Value = 1981.3846
$ws.Range("J132").Value = 2341.8333
$ws.Range("K132").Value = 5944.1538
$ws.Range("L132").Value = 7025.499899999999
$ws.Range("M132").Value = -3414.1538
$ws.Range("N132").Value = -12085.4999

$ws.Range("H134").Value = 297453.75
$ws.Range("I134").Value = 3477.7407
$ws.Range("J134").Value = 1431361.1
$ws.Range("K134").Value = 10433.2221
$ws.Range("L134").Value = 4294083.300000001
$ws.Range("M134").Value = -7898.222099999999
$ws.Range("N134").Value = -4299153.300000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1323.3143
$ws.Range("I5").Value = 404.66666
$ws.Range("J5").Value = 1802.6086
$ws.Range("K5").Value = 1213.99998
$ws.Range("L5").Value = 5407.825800000001
$ws.Range("M5").Value = -1101.99998
$ws.Range("N5").Value = -5631.825800000001

$ws.Range("H12").Value = 3366800
$ws.Range("J12").Value = 52750.21
$ws.Range("L12").Value = 158250.63
$ws.Range("N12").Value = -158596.63

$ws.Range("H131").Value = 1923987.4
$ws.Range("J131").Value = 1109.7941
$ws.Range("L131").Value = 3329.3823
$ws.Range("N131").Value = -13409.3823

$ws.Range("H135").Value = 1323.3143
$ws.Range("I135").Value = 404.66666
$ws.Range("J135").Value = 1802.6086
$ws.Range("K135").Value = 3641.99994
$ws.Range("L135").Value = 16223.4774
$ws.Range("M135").Value = -1106.99994
$ws.Range("N135").Value = -21293.4774

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6922.0454
$ws.Range("I80").Value = 9448.929
$ws.Range("K80").Value = 9448.929
$ws.Range("M80").Value = -8450.929

$ws.Range("H83").Value = 6922.0454
$ws.Range("I83").Value = 9448.929
$ws.Range("K83").Value = 47244.645
$ws.Range("M83").Value = -42252.645

$ws.Range("H113").Value = 41668220
$ws.Range("I113").Value = 76924160
$ws.Range("J113").Value = 2110.2727
$ws.Range("K113").Value = 76924160
$ws.Range("L113").Value = 2110.2727
$ws.Range("M113").Value = -76921990
$ws.Range("N113").Value = -6450.2727

$ws.Range("H132").Value = 2518.6233
$ws.Range("I132").Value = 2489.2666
$ws.Range("J132").Value = 2573.6667
$ws.Range("K132").Value = 7467.7998
$ws.Range("L132").Value = 7721.000100000001
$ws.Range("M132").Value = -4937.7998
$ws.Range("N132").Value = -12781.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 49891.43
$ws.Range("I7").Value = 73428.92999999999
$ws.Range("J7").Value = 2816.4285
$ws.Range("K7").Value = 73428.92999999999
$ws.Range("L7").Value = 2816.4285
$ws.Range("M7").Value = -73316.92999999999
$ws.Range("N7").Value = -3040.4285

$ws.Range("H18").Value = 6666.6665
$ws.Range("I18").Value = 6000
$ws.Range("K18").Value = 6000
$ws.Range("M18").Value = -5828

$ws.Range("H40").Value = 58826404
$ws.Range("I40").Value = 66669550
$ws.Range("J40").Value = 2800
$ws.Range("K40").Value = 66669550
$ws.Range("L40").Value = 2800
$ws.Range("M40").Value = -66669414
$ws.Range("N40").Value = -3072

$ws.Range("H82").Value = 796638.9399999999
$ws.Range("I82").Value = 2002648.4
$ws.Range("J82").Value = 126633.664
$ws.Range("K82").Value = 2002648.4
$ws.Range("L82").Value = 126633.664
$ws.Range("M82").Value = -2002287.4
$ws.Range("N82").Value = -127355.664

$ws.Range("H85").Value = 796638.9399999999
$ws.Range("I85").Value = 2002648.4
$ws.Range("J85").Value = 126633.664
$ws.Range("K85").Value = 2002648.4
$ws.Range("L85").Value = 126633.664
$ws.Range("M85").Value = -2001400.4
$ws.Range("N85").Value = -129129.664

$ws.Range("H122").Value = 1854593.4
$ws.Range("I122").Value = 2467759.5
$ws.Range("J122").Value = 669138.6
$ws.Range("K122").Value = 7403278.5
$ws.Range("L122").Value = 2007415.8
$ws.Range("M122").Value = -7400828.5
$ws.Range("N122").Value = -2012315.8

$ws.Range("H126").Value = 49891.43
$ws.Range("I126").Value = 73428.92999999999
$ws.Range("J126").Value = 2816.4285
$ws.Range("K126").Value = 220286.79
$ws.Range("L126").Value = 8449.2855
$ws.Range("M126").Value = -217816.79
$ws.Range("N126").Value = -13389.2855

$ws.Range("H132").Value = 8552324
$ws.Range("I132").Value = 11116868
$ws.Range("J132").Value = 3844.3333
$ws.Range("K132").Value = 33350604
$ws.Range("L132").Value = 11532.9999
$ws.Range("M132").Value = -33348074
$ws.Range("N132").Value = -16592.9999

$ws.Range("H136").Value = 7805.0454
$ws.Range("I136").Value = 6354.9287
$ws.Range("J136").Value = 10342.75
$ws.Range("K136").Value = 19064.7861
$ws.Range("L136").Value = 31028.25
$ws.Range("M136").Value = -16514.7861
$ws.Range("N136").Value = -36128.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1304.6
$ws.Range("I113").Value = 1017.5161
$ws.Range("J113").Value = 1940.2858
$ws.Range("K113").Value = 3052.5483
$ws.Range("L113").Value = 5820.857400000001
$ws.Range("M113").Value = -882.5483000000004
$ws.Range("N113").Value = -10160.8574

$ws.Range("H126").Value = 593.3
$ws.Range("I126").Value = 473.46155
$ws.Range("J126").Value = 1372.25
$ws.Range("K126").Value = 1420.38465
$ws.Range("L126").Value = 4116.75
$ws.Range("M126").Value = 1049.61535
$ws.Range("N126").Value = -9056.75

$ws.Range("H132").Value = 1561.579
$ws.Range("I132").Value = 934.56665
$ws.Range("J132").Value = 3912.875
$ws.Range("K132").Value = 2803.69995
$ws.Range("L132").Value = 11738.625
$ws.Range("M132").Value = -273.6999500000002
$ws.Range("N132").Value = -16798.625

$ws.Range("H136").Value = 2022.1719
$ws.Range("I136").Value = 2093.9211
$ws.Range("J136").Value = 1917.3077
$ws.Range("K136").Value = 6281.763300000001
$ws.Range("L136").Value = 5751.9231
$ws.Range("M136").Value = -3731.763300000001
$ws.Range("N136").Value = -10851.9231
